$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.126.52"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "2.945.51"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Formula = "'374.34"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").Formula = "'101.28"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7
$ws.Range("E7").Value = "  -0.97%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Formula = "'0.585"
$ws.Range("E9").Value = "  -1.58%  "

# Row 10
$ws.Range("D10").Formula = "'36.31"
$ws.Range("E10").Value = "  -2.52%  "

# Row 12
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$ws.Range("D13").Value = "3.407.35"
$ws.Range("E13").Value = "  -1.27%  "

# Row 14
$ws.Range("D14").Formula = "'18.01"
$ws.Range("E14").Value = "  -2.27%  "

# Row 15
$ws.Range("D15").Formula = "'7.56"
$ws.Range("E15").Value = "  -0.51%  "

# Row 16
$ws.Range("D16").Value = "2.930.04"
$ws.Range("E16").Value = "  -1.73%  "

# Row 17
$ws.Range("D17").Formula = "'11.04"
$ws.Range("E17").Value = "  +48.48%  "

# Row 18
$ws.Range("D18").Formula = "'0.979"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("D19").Value = "51.070.87"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("E20").Value = "  -5.48%  "

# Row 21
$ws.Range("D21").Formula = "'12.43"
$ws.Range("E21").Value = "  -4.13%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  -0.63%  "

# Row 23
$ws.Range("D23").Formula = "'265.25"
$ws.Range("E23").Value = "  +1.15%  "

# Row 24
$ws.Range("D24").Formula = "'68.74"
$ws.Range("E24").Value = "  -0.90%  "

# Row 25
$ws.Range("D25").Formula = "'3.14"
$ws.Range("E25").Value = "  +10.49%  "

# Row 26
$ws.Range("D26").Formula = "'8.15"
$ws.Range("E26").Value = "  -1.08%  "

# Row 27
$ws.Range("D27").Formula = "'7.53"
$ws.Range("E27").Value = "  -2.82%  "

# Row 28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Formula = "'0.111"
$ws.Range("E29").Value = "  -4.50%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Formula = "'25.68"
$ws.Range("E30").Value = "  -0.99%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Formula = "'0.164"
$ws.Range("E31").Value = "  -4.13%  "

# Row 32
$ws.Range("D32").Formula = "'10.01"
$ws.Range("E32").Value = "  +1.13%  "

# Row 33
$ws.Range("D33").Formula = "'50.69"
$ws.Range("E33").Value = "  -0.70%  "

# Row 34
$ws.Range("E34").Value = "  -1.80%  "

# Row 35
$ws.Range("D35").Formula = "'33.38"
$ws.Range("E35").Value = "  -4.96%  "

# Row 36
$ws.Range("D36").Formula = "'0.0443"
$ws.Range("E36").Value = "  -0.55%  "

# Row 37
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").Formula = "'3.17"
$ws.Range("E38").Value = "  +4.09%  "

# Row 39
$ws.Range("E39").Value = "  -0.22%  "

# Row 40
$ws.Range("D40").Formula = "'16.34"
$ws.Range("E40").Value = "  -5.12%  "

# Row 41
$ws.Range("D41").Formula = "'2.49"
$ws.Range("E41").Value = "  -4.19%  "

# Row 42
$ws.Range("E42").Value = "  -3.29%  "

# Row 43
$ws.Range("D43").Formula = "'119.98"
$ws.Range("E43").Value = "  -4.36%  "

# Row 44
$ws.Range("D44").Formula = "'21.33"
$ws.Range("E44").Value = "  -1.03%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Formula = "'3.36"
$ws.Range("E45").Value = "  +3.54%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Formula = "'0.274"
$ws.Range("E46").Value = "  -3.09%  "

# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Formula = "'2.03"
$ws.Range("E47").Value = "  -1.71%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.995.62"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Formula = "'2.29"
$ws.Range("E49").Value = "  -3.51%  "

# Row 50
$ws.Range("D50").Formula = "'0.0326"
$ws.Range("E50").Value = "  -2.64%  "

# Row 51
$ws.Range("D51").Formula = "'1.31"
$ws.Range("E51").Value = "  +1.43%  "
